# Applies the "update tech doc with more dataset references" edit:
#   1. Bumps the DATE: line from 2024-06-18 to 2024-07-02.
#   2. Inserts several new bibliography-style "Data Sources" entries,
#      each as its own BodyText-styled paragraph, in the same relative
#      position the canonical diff places them.

$d = $word.ActiveDocument

function Insert-ParaBefore($anchorText, [string[]]$newLines) {
    $r = $d.Content
    $found = $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $anchorText"
    }
    $pr = $r.Paragraphs(1).Range
    $pr.Collapse(1)
    $text = ($newLines -join "`r") + "`r"
    $pr.InsertBefore($text)
}

# 1. Update the DATE field. The date run sits right after a lone-space
#    run with identical (empty) formatting, so a Find/Replace on this
#    engine would coalesce the two into one run; clearing the found
#    range and re-inserting the new text instead keeps the original
#    three-run split ("DATE:" / " " / date) intact, matching the diff.
$r = $d.Content
$r.Find.Execute("2024-06-18", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Text = ""
$r.InsertAfter("2024-07-02")

# 2. New "Data Sources" reference paragraphs, inserted before the first
#    alphabetically-later entry that already existed (matching the diff).

Insert-ParaBefore "BC Stats. [creator] (2023). BC Demographic Survey. E01." @(
    "BC Housing. [creator] (2023). Private Market Rent Supplements. E05. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024).",
    "BC Housing. [creator] (2024). Supportive Housing. E02. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
)

Insert-ParaBefore "Ministry of Children and Family Development. [creator] (2022). Child Care Subsidy. E01." @(
    "Human Early Learning Partnership. [creator] (2024). Early Development Instrument. E03. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
)

Insert-ParaBefore "Ministry of Education and Child Care. [creator] (2023). K to 12 Socio Economic Status Index. E03." @(
    "Ministry of Education and Child Care. [creator] (2023). Affordable Child Care Benefit. E02. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
)

Insert-ParaBefore "Ministry of Health. [creator] (2022). BC Vital Events and Statistics. E04." @(
    "Ministry of Finance. [creator] (2024). Neighbourhood Income. E01. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
)

Insert-ParaBefore "Ministry of Health. [creator] (2023). Hospital Discharges. E01." @(
    "Ministry of Health. [creator] (2024). Home and Community Care. E03. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
)

Insert-ParaBefore "Ministry of Health. [creator] (2019). Registration and Premium Billings. E02." @(
    "Ministry of Health. [creator] (2020). PharmaCare. E01. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024).",
    "Ministry of Health. [creator] (2023). PharmaNet. E02. Data Innovation Program, Province of British Columbia [publisher]. Data Extract. Approver Year (2024)."
)

Write-Output "edit complete"
